$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.770.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.99%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.078.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.61%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.59"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.393"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0786"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.46%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.85"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.02"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.774"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.043.12"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "37.720.06"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.94%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.62"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.39"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.41"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.70"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.21"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.135"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.07%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.48"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.69"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0636"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.67"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.83"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.40"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.30%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.36"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0984"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.62"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.24%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0215"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.65"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.442.25"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.74%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.41"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.98"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.271.51"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.59"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.10%  "

